# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns with latest values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value.
$updates = [ordered]@{
    'D2' = '29.576.31'
    'E2' = '  +0.37%  '
    'D3' = '1.924.89'
    'E3' = '  +0.61%  '
    'D5' = '326.25'
    'E5' = '  +0.16%  '
    'E6' = '  +0.40%  '
    'D7' = '0.4817'
    'E7' = '  -0.01%  '
    'D8' = '0.4056'
    'E8' = '  -0.13%  '
    'D9' = '0.08228'
    'E9' = '  +1.02%  '
    'E10' = '  +0.04%  '
    'D11' = '23.86'
    'E11' = '  +1.73%  '
    'D12' = '1.910.46'
    'E12' = '  -0.97%  '
    'D13' = '6.118'
    'E13' = '  +1.98%  '
    'D14' = '7.309'
    'E14' = '  +2.54%  '
    'D15' = '91.83'
    'E15' = '  +1.88%  '
    'D16' = '0.06872'
    'E16' = '  +1.48%  '
    'E17' = '  +0.39%  '
    'E18' = '  +0.02%  '
    'D19' = '17.64'
    'E19' = '  -0.24%  '
    'E20' = '  +0.37%  '
    'D21' = '29.566.34'
    'E21' = '  +0.27%  '
    'D22' = '5.682'
    'E22' = '  +1.03%  '
    'D23' = '12.01'
    'E24' = '  +0.20%  '
    'D25' = '2.155.51'
    'E25' = '  +0.30%  '
    'D26' = '156.32'
    'E26' = '  +0.32%  '
    'D27' = '6.428'
    'E27' = '  +1.01%  '
    'D28' = '20.06'
    'E28' = '  +0.01%  '
    'E29' = '  -0.61%  '
    'D30' = '120.85'
    'E30' = '  +0.95%  '
    'E31' = '  -0.91%  '
    'D32' = '0.09604'
    'E32' = '  +0.81%  '
    'D33' = '5.614'
    'E33' = '  +1.80%  '
    'D34' = '3.566'
    'E34' = '  +0.13%  '
    'D35' = '1.380'
    'E35' = '  -0.52%  '
    'D36' = '0.06372'
    'E36' = '  +4.53%  '
    'D37' = '0.02295'
    'E37' = '  +1.36%  '
    'D38' = '1.193'
    'E38' = '  +1.40%  '
    'D39' = '0.5962'
    'E39' = '  -0.04%  '
    'D40' = '10.73'
    'E40' = '  +0.58%  '
    'D41' = '7.879'
    'E41' = '  -1.10%  '
    'D42' = '0.1852'
    'E42' = '  -0.05%  '
    'D43' = '2.425'
    'E43' = '  +1.11%  '
    'D44' = '1.280'
    'E44' = '  -0.41%  '
    'D45' = '12.45'
    'E45' = '  -0.68%  '
    'E46' = '  -1.10%  '
    'D47' = '0.5567'
    'E47' = '  -0.07%  '
    'D48' = '1.995'
    'E48' = '  +2.93%  '
    'D49' = '119.52'
    'E49' = '  +3.30%  '
    'D50' = '2.437'
    'E50' = '  +0.84%  '
    'D51' = '72.04'
    'E51' = '  -0.56%  '
}

foreach ($cellRef in $updates.Keys) {
    $c = $ws.Range($cellRef)
    # Force text format so numeric-looking strings (e.g. "326.25") are not
    # auto-converted into floating point numbers, losing their original text form,
    # matching the source data which stores these as text, not numbers.
    $c.NumberFormat = "@"
    $c.Value = $updates[$cellRef]
    # Restore the default "Normal" style so no extra formatting is left on the cell.
    $c.Style = "Normal"
}
